$wb = $excel.ActiveWorkbook

# --- Sheet "Ridge" (1): fix D2 value ---
$wsRidge = $wb.Worksheets.Item(1)
$wsRidge.Range("D2").Value = 30245804

# --- Sheet "Elastic net" (3): remove empty row 7 ---
$wsElastic = $wb.Worksheets.Item(3)
$wsElastic.Rows.Item(7).Delete()

# --- Sheet "Arboles" (4): add hyper-parameter header row (tree_depth, min_n) ---
$wsArboles = $wb.Worksheets.Item(4)
$arrArboles = New-Object 'object[,]' 1,12
$arrArboles[0,0]  = "nombre resultados"
$arrArboles[0,1]  = "tanda"
$arrArboles[0,2]  = "datos_usados"
$arrArboles[0,3]  = "MAE_comp"
$arrArboles[0,4]  = "Leaderboard"
$arrArboles[0,5]  = "type"
$arrArboles[0,6]  = "tree_depth"
$arrArboles[0,7]  = "min_n"
$arrArboles[0,8]  = "formula"
$arrArboles[0,9]  = "fecha"
$arrArboles[0,10] = "submittor"
$arrArboles[0,11] = "numero de modelo"
$rngArboles = $wsArboles.Range("A1:L1")
$rngArboles.Value = $arrArboles
$rngArboles.Font.Bold = $true
$wsArboles.Columns.Item(1).ColumnWidth = 17.1796875
$wsArboles.Columns.Item(12).ColumnWidth = 20.1796875

# --- Sheet "Forest" (5): add hyper-parameter header row (mtry, trees) ---
$wsForest = $wb.Worksheets.Item(5)
$arrForest = New-Object 'object[,]' 1,13
$arrForest[0,0]  = "nombre resultados"
$arrForest[0,1]  = "tanda"
$arrForest[0,2]  = "datos_usados"
$arrForest[0,3]  = "MAE_comp"
$arrForest[0,4]  = "Leaderboard"
$arrForest[0,5]  = "type"
$arrForest[0,6]  = "mtry"
$arrForest[0,7]  = "trees"
$arrForest[0,8]  = "min_n"
$arrForest[0,9]  = "formula"
$arrForest[0,10] = "fecha"
$arrForest[0,11] = "submittor"
$arrForest[0,12] = "numero de modelo"
$rngForest = $wsForest.Range("A1:M1")
$rngForest.Value = $arrForest
$rngForest.Font.Bold = $true
$wsForest.Columns.Item(13).ColumnWidth = 20.1796875

# --- Sheet "Boosting" (6): add hyper-parameter header row (learn_rate, trees, min_n) ---
$wsBoosting = $wb.Worksheets.Item(6)
$arrBoosting = New-Object 'object[,]' 1,13
$arrBoosting[0,0]  = "nombre resultados"
$arrBoosting[0,1]  = "tanda"
$arrBoosting[0,2]  = "datos_usados"
$arrBoosting[0,3]  = "MAE_comp"
$arrBoosting[0,4]  = "Leaderboard"
$arrBoosting[0,5]  = "type"
$arrBoosting[0,6]  = "learn_rate"
$arrBoosting[0,7]  = "trees"
$arrBoosting[0,8]  = "min_n"
$arrBoosting[0,9]  = "formula"
$arrBoosting[0,10] = "fecha"
$arrBoosting[0,11] = "submittor"
$arrBoosting[0,12] = "numero de modelo"
$rngBoosting = $wsBoosting.Range("A1:M1")
$rngBoosting.Value = $arrBoosting
$rngBoosting.Font.Bold = $true
$wsBoosting.Columns.Item(13).ColumnWidth = 21.81640625

# --- View state: Ridge loses the active-tab selection, Boosting becomes the active/selected sheet ---
$wsRidge.Activate()
$wsRidge.Range("A1:XFD1").Select()

$wsBoosting.Activate()
$wsBoosting.Range("K11").Select()
